$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --------------------------------------------------------------------------
# 1) Two pairs of rows in the source data had been scraped/sorted in the
#    opposite order from the canonical ordering. Swap the match-detail
#    columns (F:V) between each pair; columns A-E (index/country/league/
#    season/date) are identical within each pair so they are left as-is.
# --------------------------------------------------------------------------

# Swap F:V content between row 26 and row 27
$ws.Range("F26").Value = 'Union La Calera'
$ws.Range("G26").Value = 3
$ws.Range("H26").Value = 'Curico Unido'
$ws.Range("I26").Value = 2
$ws.Range("J26").Value = 2.34
$ws.Range("K26").Value = '04/02/2023 22:12'
$ws.Range("L26").Value = 2.04
$ws.Range("M26").Value = '11/02/2023 00:57'
$ws.Range("N26").Value = 3.62
$ws.Range("O26").Value = '04/02/2023 22:12'
$ws.Range("P26").Value = 3.42
$ws.Range("Q26").Value = '11/02/2023 00:58'
$ws.Range("R26").Value = 3.03
$ws.Range("S26").Value = '04/02/2023 22:12'
$ws.Range("T26").Value = 3.99
$ws.Range("U26").Value = '11/02/2023 00:57'
$ws.Range("V26").Value = 'https://www.betexplorer.com/football/chile/primera-division/union-la-calera-curico-unido/4Mz9Ngnb/'
$ws.Range("F27").Value = 'Copiapo'
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 'Palestino'
$ws.Range("I27").Value = 3
$ws.Range("J27").Value = 3.43
$ws.Range("K27").Value = '06/02/2023 00:42'
$ws.Range("L27").Value = 3.18
$ws.Range("M27").Value = '11/02/2023 00:53'
$ws.Range("N27").Value = 3.91
$ws.Range("O27").Value = '06/02/2023 00:42'
$ws.Range("P27").Value = 3.62
$ws.Range("Q27").Value = '11/02/2023 00:53'
$ws.Range("R27").Value = 1.94
$ws.Range("S27").Value = '06/02/2023 00:42'
$ws.Range("T27").Value = 2.26
$ws.Range("U27").Value = '11/02/2023 00:53'
$ws.Range("V27").Value = 'https://www.betexplorer.com/football/chile/primera-division/copiapo-palestino/WGWP0Wvo/'

# Swap F:V content between row 45 and row 46
$ws.Range("F45").Value = 'Everton'
$ws.Range("G45").Value = 3
$ws.Range("H45").Value = 'A. Italiano'
$ws.Range("I45").Value = 1
$ws.Range("J45").Value = 1.87
$ws.Range("K45").Value = '20/02/2023 00:41'
$ws.Range("L45").Value = 1.76
$ws.Range("M45").Value = '28/02/2023 00:28'
$ws.Range("N45").Value = 3.62
$ws.Range("O45").Value = '20/02/2023 00:41'
$ws.Range("P45").Value = 3.73
$ws.Range("Q45").Value = '28/02/2023 00:29'
$ws.Range("R45").Value = 4.48
$ws.Range("S45").Value = '20/02/2023 00:41'
$ws.Range("T45").Value = 4.99
$ws.Range("U45").Value = '28/02/2023 00:29'
$ws.Range("V45").Value = 'https://www.betexplorer.com/football/chile/primera-division/everton-a-italiano/EifizXV3/'
$ws.Range("F46").Value = 'Copiapo'
$ws.Range("G46").Value = 2
$ws.Range("H46").Value = 'O''Higgins'
$ws.Range("I46").Value = 2
$ws.Range("J46").Value = 3.21
$ws.Range("K46").Value = '20/02/2023 16:42'
$ws.Range("L46").Value = 2.51
$ws.Range("M46").Value = '28/02/2023 00:05'
$ws.Range("N46").Value = 3.48
$ws.Range("O46").Value = '20/02/2023 16:42'
$ws.Range("P46").Value = 3.32
$ws.Range("Q46").Value = '28/02/2023 00:11'
$ws.Range("R46").Value = 2.16
$ws.Range("S46").Value = '20/02/2023 16:42'
$ws.Range("T46").Value = 3
$ws.Range("U46").Value = '28/02/2023 00:11'
$ws.Range("V46").Value = 'https://www.betexplorer.com/football/chile/primera-division/copiapo-o-higgins/08d3XUFM/'

# Swap F:V content between row 63 and row 64
$ws.Range("F63").Value = 'O''Higgins'
$ws.Range("G63").Value = 0
$ws.Range("H63").Value = 'Coquimbo'
$ws.Range("I63").Value = 2
$ws.Range("J63").Value = 2.04
$ws.Range("K63").Value = '16/03/2023 08:42'
$ws.Range("L63").Value = 2.18
$ws.Range("M63").Value = '18/03/2023 21:57'
$ws.Range("N63").Value = 3.53
$ws.Range("O63").Value = '16/03/2023 08:42'
$ws.Range("P63").Value = 3.42
$ws.Range("Q63").Value = '18/03/2023 21:57'
$ws.Range("R63").Value = 3.47
$ws.Range("S63").Value = '16/03/2023 08:42'
$ws.Range("T63").Value = 3.55
$ws.Range("U63").Value = '18/03/2023 21:57'
$ws.Range("V63").Value = 'https://www.betexplorer.com/football/chile/primera-division/o-higgins-coquimbo/ALpOcRTR/'
$ws.Range("F64").Value = 'Cobresal'
$ws.Range("G64").Value = 3
$ws.Range("H64").Value = 'Colo Colo'
$ws.Range("I64").Value = 1
$ws.Range("J64").Value = 2.99
$ws.Range("K64").Value = '16/03/2023 08:42'
$ws.Range("L64").Value = 3.25
$ws.Range("M64").Value = '18/03/2023 21:57'
$ws.Range("N64").Value = 3.36
$ws.Range("O64").Value = '16/03/2023 08:42'
$ws.Range("P64").Value = 3.43
$ws.Range("Q64").Value = '18/03/2023 21:46'
$ws.Range("R64").Value = 2.43
$ws.Range("S64").Value = '16/03/2023 08:42'
$ws.Range("T64").Value = 2.31
$ws.Range("U64").Value = '18/03/2023 21:57'
$ws.Range("V64").Value = 'https://www.betexplorer.com/football/chile/primera-division/cobresal-colo-colo/KWiB0Tq9/'

# --------------------------------------------------------------------------
# 2) Append three new match rows scraped on 13-11-2023 (rows 213-215).
#    First clone the row-212 formatting (bold/bordered index style in column
#    A, datetime style in column E) down into the new rows so the new cells
#    pick up the same cell styles already used throughout the sheet, then
#    fill in the values.
# --------------------------------------------------------------------------

$ws.Range("A212:V212").Copy() | Out-Null
$ws.Range("A213:A215").PasteSpecial(-4122) | Out-Null
$ws.Range("E212").Copy() | Out-Null
$ws.Range("E213:E215").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Row 213
$ws.Range("A213").Value = 212
$ws.Range("B213").Value = 'chile'
$ws.Range("C213").Value = 'primera-division'
$ws.Range("D213").Value = "'2023"
$ws.Range("E213").Value = 45242.89583333334
$ws.Range("F213").Value = 'Colo Colo'
$ws.Range("G213").Value = 2
$ws.Range("H213").Value = 'Union La Calera'
$ws.Range("I213").Value = 0
$ws.Range("J213").Value = 1.52
$ws.Range("K213").Value = '07/11/2023 22:12'
$ws.Range("L213").Value = 1.54
$ws.Range("M213").Value = '12/11/2023 21:29'
$ws.Range("N213").Value = 4.36
$ws.Range("O213").Value = '07/11/2023 22:12'
$ws.Range("P213").Value = 4.3
$ws.Range("Q213").Value = '12/11/2023 21:29'
$ws.Range("R213").Value = 6.29
$ws.Range("S213").Value = '07/11/2023 22:12'
$ws.Range("T213").Value = 6.39
$ws.Range("U213").Value = '12/11/2023 21:29'
$ws.Range("V213").Value = 'https://www.betexplorer.com/football/chile/primera-division/colo-colo-union-la-calera/hb97FQUH/'

# Row 214
$ws.Range("A214").Value = 213
$ws.Range("B214").Value = 'chile'
$ws.Range("C214").Value = 'primera-division'
$ws.Range("D214").Value = "'2023"
$ws.Range("E214").Value = 45243
$ws.Range("F214").Value = 'Coquimbo'
$ws.Range("G214").Value = 1
$ws.Range("H214").Value = 'Copiapo'
$ws.Range("I214").Value = 0
$ws.Range("J214").Value = 1.94
$ws.Range("K214").Value = '05/11/2023 00:42'
$ws.Range("L214").Value = 1.74
$ws.Range("M214").Value = '12/11/2023 23:55'
$ws.Range("N214").Value = 3.61
$ws.Range("O214").Value = '05/11/2023 00:42'
$ws.Range("P214").Value = 3.96
$ws.Range("Q214").Value = '12/11/2023 23:55'
$ws.Range("R214").Value = 3.97
$ws.Range("S214").Value = '05/11/2023 00:42'
$ws.Range("T214").Value = 4.82
$ws.Range("U214").Value = '12/11/2023 23:55'
$ws.Range("V214").Value = 'https://www.betexplorer.com/football/chile/primera-division/coquimbo-copiapo/ve53G6FB/'

# Row 215
$ws.Range("A215").Value = 214
$ws.Range("B215").Value = 'chile'
$ws.Range("C215").Value = 'primera-division'
$ws.Range("D215").Value = "'2023"
$ws.Range("E215").Value = 45243
$ws.Range("F215").Value = 'U. Espanola'
$ws.Range("G215").Value = 3
$ws.Range("H215").Value = 'O''Higgins'
$ws.Range("I215").Value = 3
$ws.Range("J215").Value = 2.05
$ws.Range("K215").Value = '12/11/2023 20:28'
$ws.Range("L215").Value = 2.05
$ws.Range("M215").Value = '12/11/2023 20:28'
$ws.Range("N215").Value = 3.5
$ws.Range("O215").Value = '12/11/2023 20:28'
$ws.Range("P215").Value = 3.5
$ws.Range("Q215").Value = '12/11/2023 20:28'
$ws.Range("R215").Value = 3.83
$ws.Range("S215").Value = '12/11/2023 20:28'
$ws.Range("T215").Value = 3.83
$ws.Range("U215").Value = '12/11/2023 20:28'
$ws.Range("V215").Value = 'https://www.betexplorer.com/football/chile/primera-division/u-espanola-o-higgins/zXFDYmVu/'

# --------------------------------------------------------------------------
# 3) Dimensions are recalculated automatically by the engine on save, so no
#    explicit update of the sheet <dimension> element is required here.
# --------------------------------------------------------------------------
